# Add data for 2023-12-15
# Updates 2023 (column J) violent-crime counts across the citywide summary,
# the "By Neighborhood" rollup, and the individual neighborhood sheets that
# received new incident records for that day. A couple of 2016 (column C)
# corrections are included where the source diff also touched that column.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 7338
$ws.Range('J3').Value = 7714
$ws.Range('C4').Value = 1844
$ws.Range('J4').Value = 1677
$ws.Range('J6').Value = 10542
$ws.Range('C7').Value = 28388
$ws.Range('J7').Value = 27875

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J3').Value = 53
$ws.Range('J6').Value = 275
$ws.Range('J7').Value = 421

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 463
$ws.Range('J4').Value = 93
$ws.Range('J6').Value = 648
$ws.Range('J7').Value = 1757

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 421
$ws.Range('J6').Value = 450
$ws.Range('J7').Value = 1266

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J3').Value = 144
$ws.Range('J7').Value = 398

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 258
$ws.Range('J7').Value = 852

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J6').Value = 264
$ws.Range('J7').Value = 702

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 166
$ws.Range('J7').Value = 425

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 220
$ws.Range('J6').Value = 214
$ws.Range('J7').Value = 795
$ws.Range('J8').Value = 1757
$ws.Range('J13').Value = 35
$ws.Range('J15').Value = 344
$ws.Range('J16').Value = 111
$ws.Range('J19').Value = 803
$ws.Range('J20').Value = 599
$ws.Range('J21').Value = 79
$ws.Range('J22').Value = 66
$ws.Range('J24').Value = 93
$ws.Range('J29').Value = 1487
$ws.Range('J31').Value = 290
$ws.Range('J33').Value = 1266
$ws.Range('J36').Value = 376
$ws.Range('J37').Value = 852
$ws.Range('J39').Value = 21
$ws.Range('J48').Value = 309
$ws.Range('J49').Value = 172
$ws.Range('J50').Value = 166
$ws.Range('J51').Value = 350
$ws.Range('J52').Value = 712
$ws.Range('J53').Value = 421
$ws.Range('J54').Value = 549
$ws.Range('J55').Value = 439
$ws.Range('J61').Value = 30
$ws.Range('C63').Value = 273
$ws.Range('J63').Value = 84
$ws.Range('J64').Value = 185
$ws.Range('J65').Value = 702
$ws.Range('J67').Value = 1025
$ws.Range('J73').Value = 269
$ws.Range('J75').Value = 83
$ws.Range('J76').Value = 396
$ws.Range('J77').Value = 197
$ws.Range('J78').Value = 322
$ws.Range('J79').Value = 766
$ws.Range('J80').Value = 51
$ws.Range('J84').Value = 233
$ws.Range('J85').Value = 1140
$ws.Range('J89').Value = 346
$ws.Range('J90').Value = 293
$ws.Range('J94').Value = 311
$ws.Range('J95').Value = 398
$ws.Range('J96').Value = 309
$ws.Range('J99').Value = 425
$ws.Range('C101').Value = 28388
$ws.Range('J101').Value = 27875

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J4').Value = 16
$ws.Range('J7').Value = 290

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 386
$ws.Range('J7').Value = 1025

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 70
$ws.Range('J7').Value = 233

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('J6').Value = 98
$ws.Range('J7').Value = 172

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 254
$ws.Range('J7').Value = 549

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 450
$ws.Range('J3').Value = 526
$ws.Range('J6').Value = 376
$ws.Range('J7').Value = 1487

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J3').Value = 60
$ws.Range('J7').Value = 309

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 197
$ws.Range('J6').Value = 311
$ws.Range('J7').Value = 803

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J2').Value = 72
$ws.Range('J7').Value = 396

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('J2').Value = 64
$ws.Range('J7').Value = 214

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range('J5').Value = 17
$ws.Range('J6').Value = 35

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J6').Value = 100
$ws.Range('J7').Value = 322

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J6').Value = 248
$ws.Range('J7').Value = 439

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J2').Value = 31
$ws.Range('J7').Value = 93

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('J2').Value = 92
$ws.Range('J7').Value = 309

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('J6').Value = 53
$ws.Range('J7').Value = 79

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J3').Value = 256
$ws.Range('J7').Value = 766

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J6').Value = 65
$ws.Range('J7').Value = 185

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 165
$ws.Range('J7').Value = 599

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 124
$ws.Range('J7').Value = 376

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 251
$ws.Range('J4').Value = 33
$ws.Range('J7').Value = 795

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J3').Value = 58
$ws.Range('J6').Value = 165
$ws.Range('J7').Value = 311

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J6').Value = 159
$ws.Range('J7').Value = 344

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J6').Value = 57
$ws.Range('J7').Value = 166

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('J2').Value = 5
$ws.Range('J6').Value = 21

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J6').Value = 100
$ws.Range('J7').Value = 269

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J6').Value = 80
$ws.Range('J7').Value = 220

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 104
$ws.Range('J6').Value = 106
$ws.Range('J7').Value = 346

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('J3').Value = 27
$ws.Range('J7').Value = 83

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J6').Value = 88
$ws.Range('J7').Value = 293

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J3').Value = 91
$ws.Range('J6').Value = 146
$ws.Range('J7').Value = 350

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 410
$ws.Range('J6').Value = 325
$ws.Range('J7').Value = 1140

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('J6').Value = 17
$ws.Range('J7').Value = 66

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J3').Value = 65
$ws.Range('J7').Value = 197

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('J3').Value = 12
$ws.Range('J7').Value = 51

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 168
$ws.Range('J6').Value = 308
$ws.Range('J7').Value = 712

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 30

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J6').Value = 86
$ws.Range('J7').Value = 111

